$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 1346.7142
$ws.Range("I38").Value = 261.6
$ws.Range("J38").Value = 1582.6086
$ws.Range("K38").Value = 784.8000000000001
$ws.Range("L38").Value = 4747.825800000001
$ws.Range("M38").Value = -412.8000000000001
$ws.Range("N38").Value = -5491.825800000001
# Row 58
$ws.Range("H58").Value = 1610.7142
$ws.Range("I58").Value = 68.75
$ws.Range("J58").Value = 3666.6667
$ws.Range("K58").Value = 206.25
$ws.Range("L58").Value = 11000.0001
$ws.Range("M58").Value = -56.25
$ws.Range("N58").Value = -11300.0001
# Row 82
$ws.Range("H82").Value = 952.625
$ws.Range("I82").Value = 660.1429000000001
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1980.4287
$ws.Range("L82").Value = 9000
$ws.Range("M82").Value = -1574.4287
$ws.Range("N82").Value = -9812
# Row 85
$ws.Range("H85").Value = 952.625
$ws.Range("I85").Value = 660.1429000000001
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1980.4287
$ws.Range("L85").Value = 9000
$ws.Range("M85").Value = -576.4287000000002
$ws.Range("N85").Value = -11808
# Row 137
$ws.Range("H137").Value = 2560.1
$ws.Range("I137").Value = 3166.8333
$ws.Range("J137").Value = 1650
$ws.Range("K137").Value = 9500.499899999999
$ws.Range("L137").Value = 4950
$ws.Range("M137").Value = -6950.499899999999
$ws.Range("N137").Value = -10050

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17246056
$ws.Range("I32").Value = 4504.635
$ws.Range("J32").Value = 166672830
$ws.Range("K32").Value = 4504.635
$ws.Range("L32").Value = 166672830
$ws.Range("M32").Value = -4217.635
$ws.Range("N32").Value = -166673404
# Row 61
$ws.Range("H61").Value = 27779496
$ws.Range("I61").Value = 33335054
$ws.Range("J61").Value = 1700
$ws.Range("K61").Value = 33335054
$ws.Range("L61").Value = 1700
$ws.Range("M61").Value = -33334842
$ws.Range("N61").Value = -2124
# Row 74
$ws.Range("H74").Value = 2339.125
$ws.Range("I74").Value = 2766.6667
$ws.Range("J74").Value = 2082.6
$ws.Range("K74").Value = 2766.6667
$ws.Range("L74").Value = 2082.6
$ws.Range("M74").Value = -1892.6667
$ws.Range("N74").Value = -3830.6
# Row 77
$ws.Range("H77").Value = 2339.125
$ws.Range("I77").Value = 2766.6667
$ws.Range("J77").Value = 2082.6
$ws.Range("K77").Value = 13833.3335
$ws.Range("L77").Value = 10413
$ws.Range("M77").Value = -9465.333500000001
$ws.Range("N77").Value = -19149
# Row 132
$ws.Range("H132").Value = 1201282.8
$ws.Range("I132").Value = 714.7568
$ws.Range("J132").Value = 4903034
$ws.Range("K132").Value = 2144.2704
$ws.Range("L132").Value = 14709102
$ws.Range("M132").Value = 385.7296000000001
$ws.Range("N132").Value = -14714162
# Row 136
$ws.Range("H136").Value = 27779496
$ws.Range("I136").Value = 33335054
$ws.Range("J136").Value = 1700
$ws.Range("K136").Value = 100005162
$ws.Range("L136").Value = 5100
$ws.Range("M136").Value = -100002612
$ws.Range("N136").Value = -10200

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 74
$ws.Range("H74").Value = 24199.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 24199.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24199.8
$ws.Range("N74").Value = -26071.8
# Row 77
$ws.Range("H77").Value = 24199.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 24199.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 72599.39999999999
$ws.Range("N77").Value = -81959.39999999999
# Row 81
$ws.Range("H81").Value = 12477.9
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 11642.111
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 11642.111
$ws.Range("M81").Value = -18939
$ws.Range("N81").Value = -13764.111
# Row 84
$ws.Range("H84").Value = 12477.9
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 11642.111
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 34926.333
$ws.Range("M84").Value = -54696
$ws.Range("N84").Value = -45534.333
# Row 134
$ws.Range("H134").Value = 3588645
$ws.Range("I134").Value = 1104.381
$ws.Range("J134").Value = 11122481
$ws.Range("K134").Value = 3313.143
$ws.Range("L134").Value = 33367443
$ws.Range("M134").Value = -778.143
$ws.Range("N134").Value = -33372513
# Row 139
$ws.Range("H139").Value = 42799.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 42799.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 42799.75
$ws.Range("N139").Value = -53079.75
# Row 140
$ws.Range("H140").Value = 35051.2
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 35051.2
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 35051.2
$ws.Range("N140").Value = -45411.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1635098.8
$ws.Range("I31").Value = 1793043.8
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1793043.8
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1792748.8
$ws.Range("N31").Value = -3590
# Row 34
$ws.Range("H34").Value = 1635098.8
$ws.Range("I34").Value = 1793043.8
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1793043.8
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1792841.8
$ws.Range("N34").Value = -3404
# Row 58
$ws.Range("H58").Value = 30303658
$ws.Range("I58").Value = 45455188
$ws.Range("J58").Value = 597.0909
$ws.Range("K58").Value = 45455188
$ws.Range("L58").Value = 597.0909
$ws.Range("M58").Value = -45454985
$ws.Range("N58").Value = -1003.0909
# Row 132
$ws.Range("H132").Value = 20835350
$ws.Range("I132").Value = 1491
$ws.Range("J132").Value = 66669840
$ws.Range("K132").Value = 4473
$ws.Range("L132").Value = 200009520
$ws.Range("M132").Value = -1943
$ws.Range("N132").Value = -200014580
# Row 136
$ws.Range("H136").Value = 30303658
$ws.Range("I136").Value = 45455188
$ws.Range("J136").Value = 597.0909
$ws.Range("K136").Value = 136365564
$ws.Range("L136").Value = 1791.2727
$ws.Range("M136").Value = -136363014
$ws.Range("N136").Value = -6891.2727
# Row 138
$ws.Range("H138").Value = 47011.8
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 47011.8
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 47011.8
$ws.Range("N138").Value = -57291.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 808.47
$ws.Range("I131").Value = 640
$ws.Range("J131").Value = 815.48956
$ws.Range("K131").Value = 1920
$ws.Range("L131").Value = 2446.46868
$ws.Range("M131").Value = 3120
$ws.Range("N131").Value = -12526.46868

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 2201
$ws.Range("I19").Value = 1002.5
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 1002.5
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = -714.5
$ws.Range("N19").Value = -3576
# Row 132
$ws.Range("H132").Value = 8268.625
$ws.Range("I132").Value = 1585.7142
$ws.Range("J132").Value = 13466.444
$ws.Range("K132").Value = 4757.142599999999
$ws.Range("L132").Value = 40399.33199999999
$ws.Range("M132").Value = -2227.142599999999
$ws.Range("N132").Value = -45459.33199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 35722556
$ws.Range("I132").Value = 63493270
$ws.Range("J132").Value = 17350.143
$ws.Range("K132").Value = 190479810
$ws.Range("L132").Value = 52050.429
$ws.Range("M132").Value = -190477280
$ws.Range("N132").Value = -57110.429
# Row 136
$ws.Range("H136").Value = 112784376
$ws.Range("I136").Value = 81635510
$ws.Range("J136").Value = 200001200
$ws.Range("K136").Value = 244906530
$ws.Range("L136").Value = 600003600
$ws.Range("M136").Value = -244903980
$ws.Range("N136").Value = -600008700
# Row 141
$ws.Range("H141").Value = 55860
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 55860
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 55860
$ws.Range("N141").Value = -66220

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 51540.2
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 51540.2
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 51540.2
$ws.Range("N82").Value = -52306.2
# Row 85
$ws.Range("H85").Value = 51540.2
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 51540.2
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 51540.2
$ws.Range("N85").Value = -54192.2
# Row 132
$ws.Range("H132").Value = 34692.383
$ws.Range("I132").Value = 74468.42999999999
$ws.Range("J132").Value = 6849.15
$ws.Range("K132").Value = 223405.29
$ws.Range("L132").Value = 20547.45
$ws.Range("M132").Value = -220875.29
$ws.Range("N132").Value = -25607.45
# Row 136
$ws.Range("H136").Value = 1995.5778
$ws.Range("I136").Value = 1379.2222
$ws.Range("J136").Value = 2406.4814
$ws.Range("K136").Value = 4137.6666
$ws.Range("L136").Value = 7219.4442
$ws.Range("M136").Value = -1587.6666
$ws.Range("N136").Value = -12319.4442
